$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace test-case names with Jira IDs (column A) ---
$ws.Range("A2").Value = "OPQA-539"
$ws.Range("A3").Value = "OPQA-540"
$ws.Range("A4").Value = "OPQA-541"
$ws.Range("A5").Value = "OPQA-745"
$ws.Range("A6").Value = "OPQA-746"
$ws.Range("A7").Value = "OPQA-747"
$ws.Range("A8").Value = "OPQA-542"

# --- Fix API path + method for the session-eviction test cases ---
$ws.Range("D5").Value = "/admin/access"
$ws.Range("E5").Value = "PUT"
$ws.Range("D6").Value = "/admin/access"
$ws.Range("E6").Value = "PUT"

# --- Hyperlink each Jira id to its bug tracker URL (order matches rId1..rId7) ---
$ws.Hyperlinks.Add($ws.Range("A8"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-542", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-542")
$ws.Hyperlinks.Add($ws.Range("A2"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-539", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-539")
$ws.Hyperlinks.Add($ws.Range("A3"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-540", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-540")
$ws.Hyperlinks.Add($ws.Range("A4"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-541", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-541")
$ws.Hyperlinks.Add($ws.Range("A5"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-745", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-745")
$ws.Hyperlinks.Add($ws.Range("A6"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-746", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-746")
$ws.Hyperlinks.Add($ws.Range("A7"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-747", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-747")

# Hyperlinks.Add auto-applies the built-in "Hyperlink" style; the source
# file keeps these cells unstyled, so put them back to Normal.
$ws.Range("A2:A8").Style = "Normal"

# --- Widen the TESTNAME column now that it holds Jira ids ---
$ws.Columns("A").ColumnWidth = 25.15

# --- Restore the cursor position recorded in the sheet view ---
$ws.Range("A16").Select()
